$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.792.37"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").Value = "2.570.65"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'302.24"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "'92.78"
$ws.Range("E6").Value = "  -3.76%  "
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "'0.544"
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("D10").Value = "'36.08"
$ws.Range("E10").Value = "  -2.50%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("D12").Value = "'7.68"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("E13").Value = "  +6.38%  "
$ws.Range("D14").Value = "2.596.28"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "'0.881"
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("D16").Value = "'14.23"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "42.855.36"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "0.0₃0992"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "'12.70"
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("D20").Value = "'6.62"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").Value = "'71.74"
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").Value = "'253.29"
$ws.Range("E22").Value = "  -5.13%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  -4.88%  "
$ws.Range("D25").Value = "'28.79"
$ws.Range("E25").Value = "  -2.44%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "'10.27"
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").Value = "'36.99"
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").Value = "'6.01"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").Value = "'154.16"
$ws.Range("E31").Value = "  +1.38%  "
$ws.Range("D32").Value = "'2.16"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").Value = "'3.39"
$ws.Range("E33").Value = "  -5.96%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "'2.74"
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("D35").Value = "'0.0800"
$ws.Range("E35").Value = "  -2.19%  "
$ws.Range("D36").Value = "'18.35"
$ws.Range("E36").Value = "  +7.49%  "
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").Value = "'23.44"
$ws.Range("E39").Value = "  -4.36%  "
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("D41").Value = "'3.41"
$ws.Range("E41").Value = "  -3.70%  "
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("D43").Value = "2.083.54"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("E44").Value = "  +26.71%  "
$ws.Range("D45").Value = "'0.997"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "'9.23"
$ws.Range("E46").Value = "  +1.79%  "
$ws.Range("D47").Value = "'84.62"
$ws.Range("E47").Value = "  -4.50%  "
$ws.Range("D48").Value = "'107.51"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("D49").Value = "'75.85"
$ws.Range("E49").Value = "  +9.45%  "
$ws.Range("D50").Value = "2.816.00"
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("D51").Value = "'0.191"
$ws.Range("E51").Value = "  +0.04%  "
